# Reorders the species-observation rows 19-24 and 26-30 on the "Artfynd"
# sheet: each destination row receives the full record (columns A, B, D, E,
# F, G, H, Q, R) that previously lived in a different source row. Row 25 is
# untouched. This mirrors an upstream "automatic update of files" re-sync
# where record identity (column A "Id") moved between spreadsheet rows
# while the surrounding constant columns (C, I, J, K, N, P, S, T, U, V, W,
# Y, Z, AA, AB, AD-AG, AT, AW-AY) stayed the same for every row in the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices used below: A=1 B=2 D=4 E=5 F=6 G=7 H=8 Q=17 R=18
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)

# Snapshot the current ("before") values for every row that participates in
# the reshuffle so we can freely overwrite cells without clobbering a value
# we still need to read for another destination row.
$rows = @(19, 20, 21, 22, 23, 24, 26, 27, 28, 29, 30)
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destination row -> source row (source row's pre-edit data now lands here)
$mapping = @{
    19 = 21
    20 = 22
    21 = 24
    22 = 27
    23 = 20
    24 = 30
    26 = 28
    27 = 26
    28 = 29
    29 = 19
    30 = 23
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $data = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($dest, $c).Value2 = $data[$c]
    }
}
